# Generate Report for Handback
#
# The d5289b5b-... file has been handed back (its handback xliffs are now in
# sync / complete), so:
#   - the "Status" shown for that file flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" on every sheet that tracks it
#     (Overview, zh-cn, de-de)
#   - the "Latest Handback DateTime" for that file is refreshed to the
#     handback-generation timestamp, per locale
#   - the stale "version mismatch" Error Detail message for that file is
#     cleared now that it is up to date

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-09-06 05:02:10"
# Clear the stale error message but keep the cell itself present (an
# empty-text cell, not an absent one) - mirrors how the table's
# "Error Detail" column looks for rows that have no error.
$wsZhCn.Range("P3").Value = "'"
$wsZhCn.Columns.Item(16).AutoFit()

# ---- de-de sheet --------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-09-06 05:02:19"
$wsDeDe.Range("P3").Value = "'"
$wsDeDe.Columns.Item(16).AutoFit()
